$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1.0, 1.0, 5.0, 1.0, 4.0, 4.0, 1.0, 3.0, 1.0, 1.0, 1.0, 1.0, 5.0, 1.0, 4.0, 1.0, 5.0, 4.0, 1.0, 5.0, 5.0, 2.0, 1.0, 1.0, 4.0, 1.0, 4.0, 3.0, 3.0, 4.0, 4.0, 1.0, 3.0, 5.0, 1.0, 4.0, 2.0, 5.0, 3.0, 4.0, 1.0, 4.0, 5.0, 4.0, 4.0, 2.0, 2.0, 1.0, 3.0, 5.0, 1.0, 1.0, 2.0, 1.0, 2.0, 4.0, 4.0, 1.0, 1.0, 5.0, 3.0, 4.0, 5.0, 1.0, 5.0, 2.0, 2.0, 1.0, 2.0, 2.0, 2.0, 2.0, 1.0, 3.0, 4.0, 1.0, 2.0, 5.0, 1.0, 1.0, 1.0, 1.0, 4.0, 5.0, 1.0, 3.0, 4.0, 4.0, 2.0, 3.0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}
